# Jeannine's log.xlsx - add a new TUESDAY entry (8/16) at the bottom of the
# Logs sheet, reusing the same section layout/style as the existing
# "TUESDAY ... Pickup PC" block (rows 29-30).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# Duplicate the "TUESDAY / Pickup PC" section (A29:F30) down to rows 39:40,
# which brings along the correct cell styles (fills/borders/number formats)
# for both the section-header separator row and the data row.
$ws.Range("A29:F30").Copy($ws.Range("A39:F40"))
$excel.CutCopyMode = 0

# Update the date for the new entry (8/16/2016 -> serial 42598).
$ws.Range("B40").Value2 = 42598

# Update the special instructions text for the new entry.
$ws.Range("F40").Value2 = "Pick up PC and Projector and 2 matts. LEAVE PORTABLE SCREEN IN ROOM. Key for room in CB 121A storeroom on green fob. Return all equipment to Bethune 201 storeroom. IF IT IS POURING RAIN OUTSIDE LEAVE ALL EQUIPMENT LOCKED IN ROOM - Just turn off and I will pick up tomorrow."

# The longer instructions need a taller row to display comfortably.
$ws.Rows.Item(40).RowHeight = 75

# Move the active selection to reflect where the user would continue typing.
$ws.Range("F49").Select()
